$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header column (H) - copy the header style (bold, centered, bordered)
# from the adjacent "Success %" header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Updated D/E (and F) values for the first (100-iteration) block, rows 2-11
# (refit NCDE results)
$ws.Range("D2").Value = 0.8053539528284916
$ws.Range("E2").Value = 0.8053539528284916

$ws.Range("D4").Value = 0.5556631087057899
$ws.Range("E4").Value = 0.5556631087057899

$ws.Range("D5").Value = 0.5899278599882648
$ws.Range("E5").Value = 0.5899278599882648

$ws.Range("D6").Value = 0.6477824414718456
$ws.Range("E6").Value = 0.6477824414718456

$ws.Range("D7").Value = 0.7558096906415734
$ws.Range("E7").Value = 0.2441903093584266

$ws.Range("D8").Value = 0.6128933165045688
$ws.Range("E8").Value = 0.3871066834954312

$ws.Range("D9").Value = 0.4008162021707603
$ws.Range("E9").Value = 0.5991837978292397

$ws.Range("D10").Value = 0.5642755985445802
$ws.Range("E10").Value = 0.4357244014554198

$ws.Range("D11").Value = 0.4850029689237454
$ws.Range("E11").Value = 0.5149970310762546
$ws.Range("F11").Value = 0.7848498225212097

# New "Label" column values for every data row (0 = Control, 1 = MDD)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
